$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 294-295, shifting existing rows 294-386 down to 296-388
$ws.Rows("294:295").Insert()

# Populate new row 294
$ws.Cells.Item(294,1).Value2 = 7
$ws.Cells.Item(294,2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(294,3).Value2 = "Ñuble"
$ws.Cells.Item(294,4).Value2 = 44900
$ws.Cells.Item(294,5).Value2 = 16
$ws.Cells.Item(294,6).Value2 = "Fruta"
$ws.Cells.Item(294,7).Value2 = 100101
$ws.Cells.Item(294,8).Value2 = "Berries"
$ws.Cells.Item(294,9).Value2 = 100112025
$ws.Cells.Item(294,10).Value2 = "Frutilla"
$ws.Cells.Item(294,11).Value2 = "Sin especificar"
$ws.Cells.Item(294,12).Value2 = "Primera"
$ws.Cells.Item(294,13).Value2 = 120
$ws.Cells.Item(294,14).Value2 = 6500
$ws.Cells.Item(294,15).Value2 = 7000
$ws.Cells.Item(294,16).Value2 = 6750
$ws.Cells.Item(294,17).Value2 = "`$/caja 7 kilos"
$ws.Cells.Item(294,18).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(294,19).Value2 = 964
$ws.Cells.Item(294,20).Value2 = 7

# Populate new row 295
$ws.Cells.Item(295,1).Value2 = 7
$ws.Cells.Item(295,2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(295,3).Value2 = "Ñuble"
$ws.Cells.Item(295,4).Value2 = 44900
$ws.Cells.Item(295,5).Value2 = 16
$ws.Cells.Item(295,6).Value2 = "Fruta"
$ws.Cells.Item(295,7).Value2 = 100101
$ws.Cells.Item(295,8).Value2 = "Berries"
$ws.Cells.Item(295,9).Value2 = 100112025
$ws.Cells.Item(295,10).Value2 = "Frutilla"
$ws.Cells.Item(295,11).Value2 = "Sin especificar"
$ws.Cells.Item(295,12).Value2 = "Segunda"
$ws.Cells.Item(295,13).Value2 = 60
$ws.Cells.Item(295,14).Value2 = 5000
$ws.Cells.Item(295,15).Value2 = 5000
$ws.Cells.Item(295,16).Value2 = 5000
$ws.Cells.Item(295,17).Value2 = "`$/caja 7 kilos"
$ws.Cells.Item(295,18).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(295,19).Value2 = 714
$ws.Cells.Item(295,20).Value2 = 7
